$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-09-01"

# Widen column A slightly to fit the new, longer row label (the stored
# OOXML width is 25.7109375; Excel's ColumnWidth setter only resolves to
# whole-pixel steps, so feed it the value that rounds to the closest
# attainable width).
$ws.Columns.Item(1).ColumnWidth = 24.9

# August row (row 9): label text drops the "(through 08-31)" suffix now
# that the month is complete, and the 2022 total ticks up by one.
$ws.Range("A9").Value = "August"
$ws.Range("I9").Value = 168

# Duplicate the bold/bordered "totals-row" style from A9 onto the new
# A10/A11 cells before anything else touches them, so the new rows pick
# up the same cell style (s="1") instead of minting a new one.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 11 becomes the new "Total" row (previously row 10), with updated
# year-to-date sums that include the new September data.
$ws.Range("A11").Value = "Total"
$ws.Range("B11").Value = 194
$ws.Range("C11").Value = 382
$ws.Range("D11").Value = 553
$ws.Range("E11").Value = 492
$ws.Range("F11").Value = 357
$ws.Range("G11").Value = 787
$ws.Range("H11").Value = 1076
$ws.Range("I11").Value = 1142

# Row 10 becomes the new "September (through 09-01)" data row. There is
# no 2015 figure (B10) for the partial month.
$ws.Range("A10").Value = "September (through 09-01)"
$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 6
$ws.Range("I10").Value = 3
